# Applies the "want to go" (想去人数, column F) count updates and the
# refreshed event image URL (column I) for the gh-pages data refresh
# commit (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (Exhibitions) ---
$ws1.Range("F4").Value = 1283
$ws1.Range("F7").Value = 986
$ws1.Range("F12").Value = 411
$ws1.Range("F14").Value = 935
$ws1.Range("F15").Value = 1813
$ws1.Range("F16").Value = 4080
$ws1.Range("F17").Value = 1196
$ws1.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202404/bqVr0ckI1713322134807.jpeg"
$ws1.Range("F19").Value = 2658
$ws1.Range("F21").Value = 1089
$ws1.Range("F22").Value = 3631
$ws1.Range("F23").Value = 777
$ws1.Range("F26").Value = 2337
$ws1.Range("F28").Value = 857
$ws1.Range("F29").Value = 171
$ws1.Range("F30").Value = 778
$ws1.Range("F31").Value = 216
$ws1.Range("F33").Value = 1374
$ws1.Range("F34").Value = 1971
$ws1.Range("F36").Value = 500
$ws1.Range("F37").Value = 70
$ws1.Range("F38").Value = 82
$ws1.Range("F39").Value = 596
$ws1.Range("F40").Value = 287
$ws1.Range("F41").Value = 88

# --- Sheet "本地生活" (Local life) ---
$ws3.Range("F2").Value = 448

# --- Sheet "全部类型" (All types) ---
$ws4.Range("F2").Value = 448
$ws4.Range("F4").Value = 1283
$ws4.Range("F6").Value = 986
$ws4.Range("F14").Value = 411
$ws4.Range("F15").Value = 935
$ws4.Range("F16").Value = 1813
$ws4.Range("F17").Value = 4080
$ws4.Range("F18").Value = 1196
$ws4.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202404/bqVr0ckI1713322134807.jpeg"
$ws4.Range("F21").Value = 2658
$ws4.Range("F23").Value = 1089
$ws4.Range("F24").Value = 3631
$ws4.Range("F25").Value = 777
$ws4.Range("F29").Value = 2337
$ws4.Range("F33").Value = 857
$ws4.Range("F34").Value = 171
$ws4.Range("F35").Value = 779
$ws4.Range("F36").Value = 216
$ws4.Range("F38").Value = 1374
$ws4.Range("F39").Value = 1971
$ws4.Range("F43").Value = 500
$ws4.Range("F44").Value = 70
$ws4.Range("F45").Value = 596
$ws4.Range("F46").Value = 287
$ws4.Range("F47").Value = 88

$wb.Save()
